$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44: new "KKR vs RR" match entry -----------------------------
# (written in this order so newly-introduced shared strings land on the
# same indices Excel originally produced)
$ws.Range("A44").Value = "KKR vs RR"
$ws.Range("B44").Value = "RR"
$ws.Range("D44").Value = "KKR"
$ws.Range("E44").Value = "Lynn"
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0

# --- Row 45: next match, only the fixture typed in so far -------------
$ws.Range("A45").Value = "CSK vs MI"

# Man of the match for row 44, filled in after the rest of the row
$ws.Range("C44").Value = "Aaron"

# --- Formatting to match the rest of the table (centered like every
#     other data row) ---------------------------------------------------
$ws.Range("A44:E44").HorizontalAlignment = -4108
$ws.Range("H44").HorizontalAlignment = -4108
$ws.Range("A45").HorizontalAlignment = -4108
$ws.Range("H45").HorizontalAlignment = -4108

# --- Extend the two hidden helper columns (I/J) down through the new
#     rows, same formulas as the row above ------------------------------
$ws.Range("I44").Formula = "=ISNUMBER(SEARCH(""MI"",A44))"
$ws.Range("J44").Formula = "=ISNUMBER(SEARCH(""DC"",A44))"
$ws.Range("I45").Formula = "=ISNUMBER(SEARCH(""MI"",A45))"
$ws.Range("J45").Formula = "=ISNUMBER(SEARCH(""DC"",A45))"

# --- Leave the selection where the author left it ----------------------
$ws.Range("H44").Select()
